$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 115
$ws.Range("H115").Value = 11117.685
$ws.Range("I115").Value = 660.94116
$ws.Range("J115").Value = 100000
$ws.Range("K115").Value = 1982.82348
$ws.Range("L115").Value = 300000
$ws.Range("M115").Value = -415.82348
$ws.Range("N115").Value = -303134
# Row 138
$ws.Range("H138").Value = 3062.575
$ws.Range("I138").Value = 1633.4762
$ws.Range("J138").Value = 4642.1055
$ws.Range("K138").Value = 4900.4286
$ws.Range("L138").Value = 13926.3165
$ws.Range("M138").Value = 239.5713999999998
$ws.Range("N138").Value = -24206.3165

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 10419461
$ws.Range("I32").Value = 10529101
$ws.Range("K32").Value = 10529101
$ws.Range("M32").Value = -10528814
# Row 45
$ws.Range("H45").Value = 3170.7144
$ws.Range("I45").Value = 4166.6665
$ws.Range("J45").Value = 2423.75
$ws.Range("K45").Value = 4166.6665
$ws.Range("L45").Value = 2423.75
$ws.Range("M45").Value = -3789.6665
$ws.Range("N45").Value = -3177.75
# Row 64
$ws.Range("H64").Value = 40000
$ws.Range("J64").Value = 40000
$ws.Range("L64").Value = 40000
$ws.Range("N64").Value = -40496
# Row 67
$ws.Range("H67").Value = 40000
$ws.Range("J67").Value = 40000
$ws.Range("L67").Value = 40000
$ws.Range("N67").Value = -41716
# Row 74
$ws.Range("H74").Value = 2346.6904
$ws.Range("I74").Value = 2299.2368
$ws.Range("J74").Value = 2797.5
$ws.Range("K74").Value = 2299.2368
$ws.Range("L74").Value = 2797.5
$ws.Range("M74").Value = -1425.2368
$ws.Range("N74").Value = -4545.5
# Row 77
$ws.Range("H77").Value = 2346.6904
$ws.Range("I77").Value = 2299.2368
$ws.Range("J77").Value = 2797.5
$ws.Range("K77").Value = 11496.184
$ws.Range("L77").Value = 13987.5
$ws.Range("M77").Value = -7128.184000000001
$ws.Range("N77").Value = -22723.5
# Row 110
$ws.Range("H110").Value = 902.6429000000001
$ws.Range("I110").Value = 895.1539
$ws.Range("J110").Value = 1000
$ws.Range("K110").Value = 895.1539
$ws.Range("L110").Value = 1000
$ws.Range("M110").Value = 1149.8461
$ws.Range("N110").Value = -5090
# Row 132
$ws.Range("H132").Value = 4010.7
$ws.Range("I132").Value = 5774.75
$ws.Range("J132").Value = 2834.6667
$ws.Range("K132").Value = 17324.25
$ws.Range("L132").Value = 8504.000100000001
$ws.Range("M132").Value = -14794.25
$ws.Range("N132").Value = -13564.0001

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 5336.273
$ws.Range("I105").Value = 2849
$ws.Range("K105").Value = 2849
$ws.Range("M105").Value = -1102
# Row 134
$ws.Range("H134").Value = 2760.4443
$ws.Range("I134").Value = 2760.4443
$ws.Range("K134").Value = 8281.332900000001
$ws.Range("M134").Value = -5746.332900000001

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 8077.6
$ws.Range("I16").Value = 10612.571
$ws.Range("K16").Value = 10612.571
$ws.Range("M16").Value = -10325.571
# Row 31
$ws.Range("H31").Value = 2401.0476
$ws.Range("I31").Value = 2261
$ws.Range("J31").Value = 2996.25
$ws.Range("K31").Value = 2261
$ws.Range("L31").Value = 2996.25
$ws.Range("M31").Value = -1966
$ws.Range("N31").Value = -3586.25
# Row 34
$ws.Range("H34").Value = 2401.0476
$ws.Range("I34").Value = 2261
$ws.Range("J34").Value = 2996.25
$ws.Range("K34").Value = 2261
$ws.Range("L34").Value = 2996.25
$ws.Range("M34").Value = -2059
$ws.Range("N34").Value = -3400.25
# Row 62
$ws.Range("H62").Value = 8445
$ws.Range("I62").Value = 8928
$ws.Range("J62").Value = 7599.75
$ws.Range("K62").Value = 8928
$ws.Range("L62").Value = 7599.75
$ws.Range("M62").Value = -8304
$ws.Range("N62").Value = -8847.75
# Row 65
$ws.Range("H65").Value = 8445
$ws.Range("I65").Value = 8928
$ws.Range("J65").Value = 7599.75
$ws.Range("K65").Value = 44640
$ws.Range("L65").Value = 37998.75
$ws.Range("M65").Value = -41520
$ws.Range("N65").Value = -44238.75
# Row 86
$ws.Range("H86").Value = 63708.855
$ws.Range("I86").Value = 81748.25
$ws.Range("J86").Value = 39656.332
$ws.Range("K86").Value = 81748.25
$ws.Range("L86").Value = 39656.332
$ws.Range("M86").Value = -80625.25
$ws.Range("N86").Value = -41902.332
# Row 89
$ws.Range("H89").Value = 63708.855
$ws.Range("I89").Value = 81748.25
$ws.Range("J89").Value = 39656.332
$ws.Range("K89").Value = 408741.25
$ws.Range("L89").Value = 198281.66
$ws.Range("M89").Value = -403125.25
$ws.Range("N89").Value = -209513.66
# Row 99
$ws.Range("H99").Value = 6503.3335
$ws.Range("I99").Value = 2630.75
$ws.Range("J99").Value = 9601.4
$ws.Range("K99").Value = 2630.75
$ws.Range("L99").Value = 9601.4
$ws.Range("M99").Value = -1132.75
$ws.Range("N99").Value = -12597.4
# Row 113
$ws.Range("H113").Value = 8077.6
$ws.Range("I113").Value = 10612.571
$ws.Range("K113").Value = 10612.571
$ws.Range("M113").Value = -8442.571
# Row 126
$ws.Range("H126").Value = 6503.3335
$ws.Range("I126").Value = 2630.75
$ws.Range("J126").Value = 9601.4
$ws.Range("K126").Value = 7892.25
$ws.Range("L126").Value = 28804.2
$ws.Range("M126").Value = -5422.25
$ws.Range("N126").Value = -33744.2
# Row 132
$ws.Range("H132").Value = 2923.2666
$ws.Range("I132").Value = 2950.3635
$ws.Range("J132").Value = 2848.75
$ws.Range("K132").Value = 8851.0905
$ws.Range("L132").Value = 8546.25
$ws.Range("M132").Value = -6321.0905
$ws.Range("N132").Value = -13606.25
# Row 134
$ws.Range("H134").Value = 2234.25
$ws.Range("I134").Value = 1982.8182
$ws.Range("K134").Value = 5948.4546
$ws.Range("M134").Value = -3413.4546

$ws = $wb.Worksheets.Item("CUL")
# Row 48
$ws.Range("H48").Value = 730
$ws.Range("I48").Value = 333.33334
$ws.Range("K48").Value = 1000.00002
$ws.Range("M48").Value = -750.0000200000001
# Row 113
$ws.Range("H113").Value = 1449.8572
$ws.Range("I113").Value = 1480.75
$ws.Range("J113").Value = 1437.5
$ws.Range("K113").Value = 4442.25
$ws.Range("L113").Value = 4312.5
$ws.Range("M113").Value = -2272.25
$ws.Range("N113").Value = -8652.5
# Row 131
$ws.Range("H131").Value = 2915
$ws.Range("J131").Value = 5964.3335
$ws.Range("L131").Value = 17893.0005
$ws.Range("N131").Value = -27973.0005
# Row 137
$ws.Range("H137").Value = 3014.8125
$ws.Range("I137").Value = 1677.7142
$ws.Range("K137").Value = 5033.142599999999
$ws.Range("M137").Value = 66.85740000000078

$ws = $wb.Worksheets.Item("GSM")
# Row 68
$ws.Range("H68").Value = 50000
$ws.Range("J68").Value = 50000
$ws.Range("L68").Value = 50000
$ws.Range("N68").Value = -51622
# Row 71
$ws.Range("H71").Value = 50000
$ws.Range("J71").Value = 50000
$ws.Range("L71").Value = 150000
$ws.Range("N71").Value = -158112
# Row 113
$ws.Range("H113").Value = 7099.8184
$ws.Range("I113").Value = 4420
$ws.Range("K113").Value = 4420
$ws.Range("M113").Value = -2250
# Row 132
$ws.Range("H132").Value = 2894.5625
$ws.Range("I132").Value = 2174.8462
$ws.Range("J132").Value = 6013.3335
$ws.Range("K132").Value = 6524.5386
$ws.Range("L132").Value = 18040.0005
$ws.Range("M132").Value = -3994.5386
$ws.Range("N132").Value = -23100.0005

$ws = $wb.Worksheets.Item("LTW")
# Row 12
$ws.Range("H12").Value = 200
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 200
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 200
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -540
# Row 22
$ws.Range("H22").Value = 2092.2856
$ws.Range("I22").Value = 1875
$ws.Range("K22").Value = 1875
$ws.Range("M22").Value = -1580
# Row 27
$ws.Range("H27").Value = 2092.2856
$ws.Range("I27").Value = 1875
$ws.Range("K27").Value = 1875
$ws.Range("M27").Value = -1768
# Row 40
$ws.Range("H40").Value = 9325.286
$ws.Range("I40").Value = 9325.286
$ws.Range("K40").Value = 9325.286
$ws.Range("M40").Value = -9189.286
# Row 61
$ws.Range("H61").Value = 5314.1665
$ws.Range("J61").Value = 5187.2
$ws.Range("L61").Value = 5187.2
$ws.Range("N61").Value = -5591.2
# Row 113
$ws.Range("H113").Value = 5314.1665
$ws.Range("J113").Value = 5187.2
$ws.Range("L113").Value = 5187.2
$ws.Range("N113").Value = -9527.200000000001
# Row 122
$ws.Range("H122").Value = 6583.9443
$ws.Range("J122").Value = 7120.2
$ws.Range("L122").Value = 21360.6
$ws.Range("N122").Value = -26260.6
# Row 136
$ws.Range("H136").Value = 5259.278
$ws.Range("I136").Value = 5010.4375
$ws.Range("K136").Value = 15031.3125
$ws.Range("M136").Value = -12481.3125

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 1943.7142
$ws.Range("I126").Value = 1434.3334
$ws.Range("K126").Value = 4303.0002
$ws.Range("M126").Value = -1833.0002
# Row 132
$ws.Range("H132").Value = 3105.2666
$ws.Range("I132").Value = 2964
$ws.Range("J132").Value = 3493.75
$ws.Range("K132").Value = 8892
$ws.Range("L132").Value = 10481.25
$ws.Range("M132").Value = -6362
$ws.Range("N132").Value = -15541.25
# Row 136
$ws.Range("H136").Value = 1460.9131
$ws.Range("I136").Value = 937.1842
$ws.Range("J136").Value = 3948.625
$ws.Range("K136").Value = 2811.5526
$ws.Range("L136").Value = 11845.875
$ws.Range("M136").Value = -261.5526
$ws.Range("N136").Value = -16945.875
